$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 3974
$ws.Range("C3").Value = 4005
$ws.Range("C4").Value = 4005
$ws.Range("C5").Value = 4005
$ws.Range("C6").Value = 4081
$ws.Range("C7").Value = 4081
$ws.Range("C8").Value = 4081
$ws.Range("C9").Value = 4212
$ws.Range("C10").Value = 4212
$ws.Range("C11").Value = 4488
$ws.Range("C12").Value = 4488
